$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import Priorities")

# SubstanceInFuelMix now running with the CO2 Market again:
# add the PowerGeneratingTechnologyFuel import priority row.
$ws.Range("A11").Value = "PowerGeneratingTechnologyFuel"
$ws.Range("B11").Value = 1

# Column A needs to be a bit wider to fit the new (longer) class name.
$ws.Columns.Item(1).ColumnWidth = 29.5

# Leave the selection on the new entry's priority column, one row below it.
$ws.Range("C11").Select()
